$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values: condition_source_value becomes text "DOA", condition_source_concept_id becomes 0
$ws.Range("A2").Value = "DOA"
$ws.Range("B2").Value = 0

# Column widths / row height adjustments
$ws.Range("A1").EntireRow.RowHeight = 16.5
$ws.Range("A2").EntireRow.RowHeight = 16.5
$ws.Columns.Item(1).ColumnWidth = 22.375
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 20.5

# Selection moves to A3
$ws.Range("A3").Select()
